$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.917497143191845
$ws.Range("D2").Value = 8.379534470898294
$ws.Range("E2").Value = 13.17150488204265
$ws.Range("F2").Value = 38.31587132213524
$ws.Range("G2").Value = 3.637461528759323
$ws.Range("J2").Value = 9.9707373658851
$ws.Range("O2").Value = 29.15539173534993
$ws.Range("C3").Value = 4.750038287576726
$ws.Range("D3").Value = 8.39465060956497
$ws.Range("E3").Value = 13.15978120347594
$ws.Range("F3").Value = 37.81025459191149
$ws.Range("G3").Value = 3.641781056497856
$ws.Range("J3").Value = 9.975187776357545
$ws.Range("O3").Value = 28.79417684634626
$ws.Range("C4").Value = 4.645915767838323
$ws.Range("D4").Value = 8.404845278276241
$ws.Range("E4").Value = 13.15545742062733
$ws.Range("F4").Value = 37.50866947976637
$ws.Range("G4").Value = 3.644570570457482
$ws.Range("J4").Value = 9.979858037231621
$ws.Range("O4").Value = 28.57929408398374
$ws.Range("C5").Value = 4.603233074408466
$ws.Range("D5").Value = 8.409229529079532
$ws.Range("E5").Value = 13.15441859798488
$ws.Range("F5").Value = 37.38814548990334
$ws.Range("G5").Value = 3.645741988541843
$ws.Range("J5").Value = 9.982247850571039
$ws.Range("O5").Value = 28.49356298653688
$ws.Range("C6").Value = 4.596132744288062
$ws.Range("D6").Value = 8.40997141774162
$ws.Range("E6").Value = 13.15428977878576
$ws.Range("F6").Value = 37.36828010956901
$ws.Range("G6").Value = 3.645938599671351
$ws.Range("J6").Value = 9.982674048260337
$ws.Range("O6").Value = 28.47944106643437
$ws.Range("C7").Value = 4.645341050299908
$ws.Range("D7").Value = 8.404903474950212
$ws.Range("E7").Value = 13.15544048259849
$ws.Range("F7").Value = 37.50703424965018
$ws.Range("G7").Value = 3.644586228050325
$ws.Range("J7").Value = 9.979888297599421
$ws.Range("O7").Value = 28.57813032523613
$ws.Range("C8").Value = 4.860071148188004
$ws.Range("D8").Value = 8.384557071697474
$ws.Range("E8").Value = 13.16686625684076
$ws.Range("F8").Value = 38.1397832623194
$ws.Range("G8").Value = 3.638922489430194
$ws.Range("J8").Value = 9.971869385334447
$ws.Range("O8").Value = 29.02947253667042
$ws.Range("C9").Value = 5.267797834203716
$ws.Range("D9").Value = 8.351896488628002
$ws.Range("E9").Value = 13.21204874189657
$ws.Range("F9").Value = 39.44433865126175
$ws.Range("G9").Value = 3.628898923999385
$ws.Range("J9").Value = 9.97154006306204
$ws.Range("O9").Value = 29.96479493695715
$ws.Range("C10").Value = 5.555577424452244
$ws.Range("D10").Value = 8.33230312347232
$ws.Range("E10").Value = 13.25905958193102
$ws.Range("F10").Value = 40.43258108286882
$ws.Range("G10").Value = 3.622185812989635
$ws.Range("J10").Value = 9.980703773332717
$ws.Range("O10").Value = 30.67632919590963
$ws.Range("C11").Value = 5.683289480366705
$ws.Range("D11").Value = 8.324343274162562
$ws.Range("E11").Value = 13.28342094219021
$ws.Range("F11").Value = 40.88667420338388
$ws.Range("G11").Value = 3.619271314870923
$ws.Range("J11").Value = 9.986915426135704
$ws.Range("O11").Value = 31.00394998365647
$ws.Range("C12").Value = 5.731144146241998
$ws.Range("D12").Value = 8.321465978278935
$ws.Range("E12").Value = 13.29307089710904
$ws.Range("F12").Value = 41.05911770167091
$ws.Range("G12").Value = 3.618187556162465
$ws.Range("J12").Value = 9.989561056539991
$ws.Range("O12").Value = 31.12846399087006
$ws.Range("C13").Value = 5.720861031978075
$ws.Range("D13").Value = 8.322079568073088
$ws.Range("E13").Value = 13.29097376640005
$ws.Range("F13").Value = 41.02195977649668
$ws.Range("G13").Value = 3.618420080178943
$ws.Range("J13").Value = 9.988978231202932
$ws.Range("O13").Value = 31.10162943268644
$ws.Range("C14").Value = 5.687236925558989
$ws.Range("D14").Value = 8.324103813980873
$ws.Range("E14").Value = 13.28420635488909
$ws.Range("F14").Value = 40.90085237832001
$ws.Range("G14").Value = 3.619181755378191
$ws.Range("J14").Value = 9.987127206179997
$ws.Range("O14").Value = 31.01418544175409
$ws.Range("C15").Value = 5.666573847298634
$ws.Range("D15").Value = 8.325361550291001
$ws.Range("E15").Value = 13.28011634602458
$ws.Range("F15").Value = 40.82672936358064
$ws.Range("G15").Value = 3.619650890739402
$ws.Range("J15").Value = 9.986031596312412
$ws.Range("O15").Value = 30.96067881572742
$ws.Range("C16").Value = 5.54716260824822
$ws.Range("D16").Value = 8.332842487610503
$ws.Range("E16").Value = 13.25752713912081
$ws.Range("F16").Value = 40.40298223827926
$ws.Range("G16").Value = 3.622379074251351
$ws.Range("J16").Value = 9.980338910971648
$ws.Range("O16").Value = 30.65498779838077
$ws.Range("C17").Value = 5.473052887979478
$ws.Range("D17").Value = 8.337675848306244
$ws.Range("E17").Value = 13.24442954544496
$ws.Range("F17").Value = 40.144071592475
$ws.Range("G17").Value = 3.624088314602398
$ws.Range("J17").Value = 9.977369701394309
$ws.Range("O17").Value = 30.46838297396408
$ws.Range("C18").Value = 5.430128156078242
$ws.Range("D18").Value = 8.340545605306088
$ws.Range("E18").Value = 13.23717647786321
$ws.Range("F18").Value = 39.99558975782939
$ws.Range("G18").Value = 3.62508454699288
$ws.Range("J18").Value = 9.975854225933658
$ws.Range("O18").Value = 30.36143056204235
$ws.Range("C19").Value = 5.415544782208891
$ws.Range("D19").Value = 8.341532672509217
$ws.Range("E19").Value = 13.23476894553338
$ws.Range("F19").Value = 39.94539640463731
$ws.Range("G19").Value = 3.625424112149467
$ws.Range("J19").Value = 9.975374153825003
$ws.Range("O19").Value = 30.32528672754191
$ws.Range("C20").Value = 5.480973273934521
$ws.Range("D20").Value = 8.337152042595159
$ws.Range("E20").Value = 13.24579481525484
$ws.Range("F20").Value = 40.17158907014781
$ws.Range("G20").Value = 3.623905005908868
$ws.Range("J20").Value = 9.977665874913981
$ws.Range("O20").Value = 30.48820912084535
$ws.Range("C21").Value = 5.697127250980085
$ws.Range("D21").Value = 8.323505528986789
$ws.Range("E21").Value = 13.28618260483255
$ws.Range("F21").Value = 40.93641259953612
$ws.Range("G21").Value = 3.618957493960103
$ws.Range("J21").Value = 9.987662937952054
$ws.Range("O21").Value = 31.03985853224044
$ws.Range("C22").Value = 5.835422455774034
$ws.Range("D22").Value = 8.315384781008179
$ws.Range("E22").Value = 13.31505285444387
$ws.Range("F22").Value = 41.43904107124746
$ws.Range("G22").Value = 3.615839929079876
$ws.Range("J22").Value = 9.995906512116214
$ws.Range("O22").Value = 31.40297030428745
$ws.Range("C23").Value = 5.761897942792263
$ws.Range("D23").Value = 8.319646008059124
$ws.Range("E23").Value = 13.29941899025573
$ws.Range("F23").Value = 41.17057863238667
$ws.Range("G23").Value = 3.617493270602819
$ws.Range("J23").Value = 9.991350479635134
$ws.Range("O23").Value = 31.20897271189405
$ws.Range("C24").Value = 5.477393454532767
$ws.Range("D24").Value = 8.3373885717213
$ws.Range("E24").Value = 13.2451767140214
$ws.Range("F24").Value = 40.1591472555916
$ws.Range("G24").Value = 3.623987837499457
$ws.Range("J24").Value = 9.977531378164741
$ws.Range("O24").Value = 30.47924468572206
$ws.Range("C25").Value = 5.159325451038388
$ws.Range("D25").Value = 8.359958167409118
$ws.Range("E25").Value = 13.19739254666036
$ws.Range("F25").Value = 39.08552773273311
$ws.Range("G25").Value = 3.631495555123247
$ws.Range("J25").Value = 9.969978108515276
$ws.Range("O25").Value = 29.70702748933627
